# Scheduled market-data refresh: update cached currentAveragePrice /
# LevePrice / LeveProfit columns (H:N) for the affected Leve rows on each
# job sheet. Generated from the upstream price-bot diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 520  # H18 (was 500)
$ws.Cells.Item(18, 10).Value = 600  # J18 (was 0)
$ws.Cells.Item(18, 12).Value = 600  # L18 (was 0)
$ws.Cells.Item(18, 14).Value = -1168  # N18: new cell
$ws.Cells.Item(19, 8).Value = 1265.6522  # H19 (was 1318.7727)
$ws.Cells.Item(19, 9).Value = 1530  # I19 (was 1673.3)
$ws.Cells.Item(19, 11).Value = 1530  # K19 (was 1673.3)
$ws.Cells.Item(19, 13).Value = -1355  # M19 (was -1498.3)
$ws.Cells.Item(32, 8).Value = 4054.75  # H32 (was 4291.3)
$ws.Cells.Item(32, 10).Value = 2078.6667  # J32 (was 1852)
$ws.Cells.Item(32, 12).Value = 2078.6667  # L32 (was 1852)
$ws.Cells.Item(32, 14).Value = -2730.6667  # N32 (was -2504)
$ws.Cells.Item(40, 8).Value = 4271.793  # H40 (was 4406.926)
$ws.Cells.Item(40, 9).Value = 2681.3333  # I40 (was 2748.1428)
$ws.Cells.Item(40, 11).Value = 2681.3333  # K40 (was 2748.1428)
$ws.Cells.Item(40, 13).Value = -2506.3333  # M40 (was -2573.1428)
$ws.Cells.Item(61, 8).Value = 0  # H61 (was 250)
$ws.Cells.Item(61, 9).Value = 0  # I61 (was 250)
$ws.Cells.Item(61, 11).Value = 0  # K61 (was 750)
$ws.Cells.Item(61, 13).Value = ""  # M61: clear (was -578)
$ws.Cells.Item(100, 8).Value = 1474.1111  # H100 (was 1406.5454)
$ws.Cells.Item(100, 9).Value = 1474.1111  # I100 (was 1346.6)
$ws.Cells.Item(100, 10).Value = 0  # J100 (was 2006)
$ws.Cells.Item(100, 11).Value = 1474.1111  # K100 (was 1346.6)
$ws.Cells.Item(100, 12).Value = 0  # L100 (was 2006)
$ws.Cells.Item(100, 13).Value = -933.1111000000001  # M100 (was -805.5999999999999)
$ws.Cells.Item(100, 14).Value = ""  # N100: clear (was -3088)
$ws.Cells.Item(113, 8).Value = 34497108  # H113 (was 35729070)
$ws.Cells.Item(113, 9).Value = 83337460  # I113 (was 76927090)
$ws.Cells.Item(113, 10).Value = 21571.883  # J113 (was 24125.533)
$ws.Cells.Item(113, 11).Value = 83337460  # K113 (was 76927090)
$ws.Cells.Item(113, 12).Value = 21571.883  # L113 (was 24125.533)
$ws.Cells.Item(113, 13).Value = -83334206  # M113 (was -76923836)
$ws.Cells.Item(113, 14).Value = -28079.883  # N113 (was -30633.533)
$ws.Cells.Item(116, 8).Value = 3639.3845  # H116 (was 3420.7)
$ws.Cells.Item(116, 9).Value = 2764  # I116 (was 3167.8333)
$ws.Cells.Item(116, 10).Value = 5040  # J116 (was 3800)
$ws.Cells.Item(116, 11).Value = 2764  # K116 (was 3167.8333)
$ws.Cells.Item(116, 12).Value = 5040  # L116 (was 3800)
$ws.Cells.Item(116, 13).Value = 678  # M116 (was 274.1667000000002)
$ws.Cells.Item(116, 14).Value = -11924  # N116 (was -10684)
$ws.Cells.Item(137, 8).Value = 2163.3076  # H137 (was 2229.84)
$ws.Cells.Item(137, 9).Value = 2092.0454  # I137 (was 2167.8572)
$ws.Cells.Item(137, 11).Value = 6276.1362  # K137 (was 6503.571599999999)
$ws.Cells.Item(137, 13).Value = -3726.1362  # M137 (was -3953.571599999999)

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 4440.44  # H45 (was 4483.7915)
$ws.Cells.Item(45, 10).Value = 6351.625  # J45 (was 6773.2856)
$ws.Cells.Item(45, 12).Value = 6351.625  # L45 (was 6773.2856)
$ws.Cells.Item(45, 14).Value = -7105.625  # N45 (was -7527.2856)
$ws.Cells.Item(102, 8).Value = 4795.8184  # H102 (was 4075.9)
$ws.Cells.Item(102, 9).Value = 2095.4375  # I102 (was 2095.5625)
$ws.Cells.Item(102, 10).Value = 11996.833  # J102 (was 11997.25)
$ws.Cells.Item(102, 11).Value = 2095.4375  # K102 (was 2095.5625)
$ws.Cells.Item(102, 12).Value = 11996.833  # L102 (was 11997.25)
$ws.Cells.Item(102, 13).Value = -473.4375  # M102 (was -473.5625)
$ws.Cells.Item(102, 14).Value = -15240.833  # N102 (was -15241.25)
$ws.Cells.Item(132, 8).Value = 1497.4828  # H132 (was 1685.0741)
$ws.Cells.Item(132, 9).Value = 1427.0714  # I132 (was 1542.4)
$ws.Cells.Item(132, 10).Value = 3469  # J132 (was 3468.5)
$ws.Cells.Item(132, 11).Value = 4281.2142  # K132 (was 4627.200000000001)
$ws.Cells.Item(132, 12).Value = 10407  # L132 (was 10405.5)
$ws.Cells.Item(132, 13).Value = -1751.2142  # M132 (was -2097.200000000001)
$ws.Cells.Item(132, 14).Value = -15467  # N132 (was -15465.5)

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2827.2  # H20 (was 2540.6206)
$ws.Cells.Item(20, 9).Value = 2694.5  # I20 (was 2414.2856)
$ws.Cells.Item(20, 10).Value = 3168.4285  # J20 (was 2872.25)
$ws.Cells.Item(20, 11).Value = 2694.5  # K20 (was 2414.2856)
$ws.Cells.Item(20, 12).Value = 3168.4285  # L20 (was 2872.25)
$ws.Cells.Item(20, 13).Value = -2447.5  # M20 (was -2167.2856)
$ws.Cells.Item(20, 14).Value = -3662.4285  # N20 (was -3366.25)
$ws.Cells.Item(94, 8).Value = 840.61536  # H94 (was 894.25)
$ws.Cells.Item(94, 9).Value = 448.0909  # I94 (was 473.2)
$ws.Cells.Item(94, 11).Value = 448.0909  # K94 (was 473.2)
$ws.Cells.Item(94, 13).Value = 2.909100000000024  # M94 (was -22.19999999999999)
$ws.Cells.Item(99, 8).Value = 2454.868  # H99 (was 2420.5186)
$ws.Cells.Item(99, 9).Value = 2840.7273  # I99 (was 2743.3044)
$ws.Cells.Item(99, 11).Value = 2840.7273  # K99 (was 2743.3044)
$ws.Cells.Item(99, 13).Value = -1342.7273  # M99 (was -1245.3044)
$ws.Cells.Item(107, 8).Value = 1712.6522  # H107 (was 1788.8695)
$ws.Cells.Item(107, 9).Value = 1471.95  # I107 (was 1559.6)
$ws.Cells.Item(107, 11).Value = 1471.95  # K107 (was 1559.6)
$ws.Cells.Item(107, 13).Value = 448.05  # M107 (was 360.4000000000001)
$ws.Cells.Item(132, 8).Value = 98176.336  # H132 (was 98389.5)
$ws.Cells.Item(132, 10).Value = 98176.336  # J132 (was 98389.5)
$ws.Cells.Item(132, 12).Value = 98176.336  # L132 (was 98389.5)
$ws.Cells.Item(132, 14).Value = -108296.336  # N132 (was -108509.5)

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3721  # H16 (was 3453.5217)
$ws.Cells.Item(16, 9).Value = 2445.6924  # I16 (was 2291.7144)
$ws.Cells.Item(16, 10).Value = 5793.375  # J16 (was 5260.778)
$ws.Cells.Item(16, 11).Value = 2445.6924  # K16 (was 2291.7144)
$ws.Cells.Item(16, 12).Value = 5793.375  # L16 (was 5260.778)
$ws.Cells.Item(16, 13).Value = -2158.6924  # M16 (was -2004.7144)
$ws.Cells.Item(16, 14).Value = -6367.375  # N16 (was -5834.778)
$ws.Cells.Item(22, 8).Value = 411.8  # H22 (was 465.375)
$ws.Cells.Item(22, 9).Value = 402  # I22 (was 460.42856)
$ws.Cells.Item(22, 11).Value = 402  # K22 (was 460.42856)
$ws.Cells.Item(22, 13).Value = -52  # M22 (was -110.42856)
$ws.Cells.Item(93, 8).Value = 20284.715  # H93 (was 17332.166)
$ws.Cells.Item(93, 10).Value = 35998.332  # J93 (was 34997.5)
$ws.Cells.Item(93, 12).Value = 35998.332  # L93 (was 34997.5)
$ws.Cells.Item(93, 14).Value = -39742.332  # N93 (was -38741.5)
$ws.Cells.Item(105, 8).Value = 3249.5  # H105 (was 2783)
$ws.Cells.Item(105, 9).Value = 2999.6667  # I105 (was 2539.8)
$ws.Cells.Item(105, 11).Value = 2999.6667  # K105 (was 2539.8)
$ws.Cells.Item(105, 13).Value = -1252.6667  # M105 (was -792.8000000000002)
$ws.Cells.Item(113, 8).Value = 3721  # H113 (was 3453.5217)
$ws.Cells.Item(113, 9).Value = 2445.6924  # I113 (was 2291.7144)
$ws.Cells.Item(113, 10).Value = 5793.375  # J113 (was 5260.778)
$ws.Cells.Item(113, 11).Value = 2445.6924  # K113 (was 2291.7144)
$ws.Cells.Item(113, 12).Value = 5793.375  # L113 (was 5260.778)
$ws.Cells.Item(113, 13).Value = -275.6923999999999  # M113 (was -121.7143999999998)
$ws.Cells.Item(113, 14).Value = -10133.375  # N113 (was -9600.778)
$ws.Cells.Item(132, 8).Value = 6408.1  # H132 (was 4321.9375)
$ws.Cells.Item(132, 9).Value = 4308.8  # I132 (was 2419.4546)
$ws.Cells.Item(132, 11).Value = 12926.4  # K132 (was 7258.3638)
$ws.Cells.Item(132, 13).Value = -10396.4  # M132 (was -4728.3638)

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 1314.1428  # H8 (was 1366.8334)
$ws.Cells.Item(8, 9).Value = 1314.1428  # I8 (was 1366.8334)
$ws.Cells.Item(8, 11).Value = 3942.4284  # K8 (was 4100.5002)
$ws.Cells.Item(8, 13).Value = -3803.4284  # M8 (was -3961.5002)
$ws.Cells.Item(64, 8).Value = 2140.6667  # H64 (was 2229)
$ws.Cells.Item(64, 9).Value = 1711  # I64 (was 1715)
$ws.Cells.Item(64, 11).Value = 5133  # K64 (was 5145)
$ws.Cells.Item(64, 13).Value = -4863  # M64 (was -4875)
$ws.Cells.Item(67, 8).Value = 2140.6667  # H67 (was 2229)
$ws.Cells.Item(67, 9).Value = 1711  # I67 (was 1715)
$ws.Cells.Item(67, 11).Value = 5133  # K67 (was 5145)
$ws.Cells.Item(67, 13).Value = -4197  # M67 (was -4209)
$ws.Cells.Item(97, 8).Value = 3549.8  # H97 (was 3049.8)
$ws.Cells.Item(97, 9).Value = 2583.3333  # I97 (was 2562.25)
$ws.Cells.Item(97, 10).Value = 4999.5  # J97 (was 5000)
$ws.Cells.Item(97, 11).Value = 7749.999899999999  # K97 (was 7686.75)
$ws.Cells.Item(97, 12).Value = 14998.5  # L97 (was 15000)
$ws.Cells.Item(97, 13).Value = -7253.999899999999  # M97 (was -7190.75)
$ws.Cells.Item(97, 14).Value = -15990.5  # N97 (was -15992)
$ws.Cells.Item(98, 8).Value = 265.0909  # H98 (was 272.1111)
$ws.Cells.Item(98, 9).Value = 235  # I98 (was 240)
$ws.Cells.Item(98, 10).Value = 268.1  # J98 (was 281.2857)
$ws.Cells.Item(98, 11).Value = 705  # K98 (was 720)
$ws.Cells.Item(98, 12).Value = 804.3000000000001  # L98 (was 843.8571000000001)
$ws.Cells.Item(98, 13).Value = 793  # M98 (was 778)
$ws.Cells.Item(98, 14).Value = -3800.3  # N98 (was -3839.8571)
$ws.Cells.Item(127, 8).Value = 7332.6665  # H127 (was 8999)
$ws.Cells.Item(127, 10).Value = 7332.6665  # J127 (was 8999)
$ws.Cells.Item(127, 12).Value = 21997.9995  # L127 (was 26997)
$ws.Cells.Item(127, 14).Value = -31917.9995  # N127 (was -36917)

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2699.5386  # H80 (was 2656.2)
$ws.Cells.Item(80, 10).Value = 2703.25  # J80 (was 2637.5)
$ws.Cells.Item(80, 12).Value = 2703.25  # L80 (was 2637.5)
$ws.Cells.Item(80, 14).Value = -4699.25  # N80 (was -4633.5)
$ws.Cells.Item(83, 8).Value = 2699.5386  # H83 (was 2656.2)
$ws.Cells.Item(83, 10).Value = 2703.25  # J83 (was 2637.5)
$ws.Cells.Item(83, 12).Value = 13516.25  # L83 (was 13187.5)
$ws.Cells.Item(83, 14).Value = -23500.25  # N83 (was -23171.5)
$ws.Cells.Item(141, 8).Value = 172768.6  # H141 (was 555444)
$ws.Cells.Item(141, 9).Value = 70000  # I141 (was 0)
$ws.Cells.Item(141, 10).Value = 198460.75  # J141 (was 555444)
$ws.Cells.Item(141, 11).Value = 70000  # K141 (was 0)
$ws.Cells.Item(141, 12).Value = 198460.75  # L141 (was 555444)
$ws.Cells.Item(141, 13).Value = -64820  # M141: new cell
$ws.Cells.Item(141, 14).Value = -208820.75  # N141 (was -565804)

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2781.5715  # H40 (was 2889.879)
$ws.Cells.Item(40, 9).Value = 1837.7727  # I40 (was 1877.7142)
$ws.Cells.Item(40, 10).Value = 4378.769  # J40 (was 4661.1665)
$ws.Cells.Item(40, 11).Value = 1837.7727  # K40 (was 1877.7142)
$ws.Cells.Item(40, 12).Value = 4378.769  # L40 (was 4661.1665)
$ws.Cells.Item(40, 13).Value = -1701.7727  # M40 (was -1741.7142)
$ws.Cells.Item(40, 14).Value = -4650.769  # N40 (was -4933.1665)
$ws.Cells.Item(61, 8).Value = 80486.53999999999  # H61 (was 74856.07000000001)
$ws.Cells.Item(61, 9).Value = 80486.53999999999  # I61 (was 74856.07000000001)
$ws.Cells.Item(61, 11).Value = 80486.53999999999  # K61 (was 74856.07000000001)
$ws.Cells.Item(61, 13).Value = -80284.53999999999  # M61 (was -74654.07000000001)
$ws.Cells.Item(93, 8).Value = 10316  # H93 (was 8963.643)
$ws.Cells.Item(93, 9).Value = 2005.2941  # I93 (was 1785.1428)
$ws.Cells.Item(93, 11).Value = 2005.2941  # K93 (was 1785.1428)
$ws.Cells.Item(93, 13).Value = -757.2941000000001  # M93 (was -537.1428000000001)
$ws.Cells.Item(113, 8).Value = 80486.53999999999  # H113 (was 74856.07000000001)
$ws.Cells.Item(113, 9).Value = 80486.53999999999  # I113 (was 74856.07000000001)
$ws.Cells.Item(113, 11).Value = 80486.53999999999  # K113 (was 74856.07000000001)
$ws.Cells.Item(113, 13).Value = -78316.53999999999  # M113 (was -72686.07000000001)
$ws.Cells.Item(122, 8).Value = 4028  # H122 (was 4572.875)
$ws.Cells.Item(122, 9).Value = 2487.2856  # I122 (was 2742.8)
$ws.Cells.Item(122, 11).Value = 7461.8568  # K122 (was 8228.400000000001)
$ws.Cells.Item(122, 13).Value = -5011.8568  # M122 (was -5778.400000000001)
$ws.Cells.Item(132, 8).Value = 11588.956  # H132 (was 12359.619)
$ws.Cells.Item(132, 9).Value = 13187.053  # I132 (was 14327.059)
$ws.Cells.Item(132, 11).Value = 39561.159  # K132 (was 42981.177)
$ws.Cells.Item(132, 13).Value = -37031.159  # M132 (was -40451.177)

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 627.3333  # H100 (was 601.7143)
$ws.Cells.Item(100, 9).Value = 554.94116  # I100 (was 535.3333)
$ws.Cells.Item(100, 10).Value = 935  # J100 (was 1000)
$ws.Cells.Item(100, 11).Value = 1109.88232  # K100 (was 1070.6666)
$ws.Cells.Item(100, 12).Value = 1870  # L100 (was 2000)
$ws.Cells.Item(100, 13).Value = -568.8823199999999  # M100 (was -529.6666)
$ws.Cells.Item(100, 14).Value = -2952  # N100 (was -3082)
$ws.Cells.Item(107, 8).Value = 1091.2  # H107 (was 1049.0625)
$ws.Cells.Item(107, 9).Value = 1106.2727  # I107 (was 1115.5454)
$ws.Cells.Item(107, 10).Value = 1049.75  # J107 (was 902.8)
$ws.Cells.Item(107, 11).Value = 3318.8181  # K107 (was 3346.6362)
$ws.Cells.Item(107, 12).Value = 3149.25  # L107 (was 2708.4)
$ws.Cells.Item(107, 13).Value = -1398.8181  # M107 (was -1426.6362)
$ws.Cells.Item(107, 14).Value = -6989.25  # N107 (was -6548.4)
$ws.Cells.Item(122, 8).Value = 3026.3704  # H122 (was 3027.8076)
$ws.Cells.Item(122, 9).Value = 1891.381  # I122 (was 1891.7142)
$ws.Cells.Item(122, 10).Value = 6998.8335  # J122 (was 7799.4)
$ws.Cells.Item(122, 11).Value = 5674.143  # K122 (was 5675.142599999999)
$ws.Cells.Item(122, 12).Value = 20996.5005  # L122 (was 23398.2)
$ws.Cells.Item(122, 13).Value = -3224.143  # M122 (was -3225.142599999999)
$ws.Cells.Item(122, 14).Value = -25896.5005  # N122 (was -28298.2)
$ws.Cells.Item(140, 8).Value = 90299.62  # H140 (was 117461.5)
$ws.Cells.Item(140, 9).Value = 75000  # I140 (was 0)
$ws.Cells.Item(140, 10).Value = 91574.586  # J140 (was 117461.5)
$ws.Cells.Item(140, 11).Value = 75000  # K140 (was 0)
$ws.Cells.Item(140, 12).Value = 91574.586  # L140 (was 117461.5)
$ws.Cells.Item(140, 13).Value = -69820  # M140: new cell
$ws.Cells.Item(140, 14).Value = -101934.586  # N140 (was -127821.5)
$ws.Cells.Item(141, 8).Value = 87076.73  # H141 (was 90597.39999999999)
$ws.Cells.Item(141, 10).Value = 87076.73  # J141 (was 90597.39999999999)
$ws.Cells.Item(141, 12).Value = 87076.73  # L141 (was 90597.39999999999)
$ws.Cells.Item(141, 14).Value = -97436.73  # N141 (was -100957.4)
